# plantilla_certificado.docx edit
#
# 1. Collapse the "+++" / "domicilio" / "+++" run-triplet (and the same
#    for "localidad") into a single run per placeholder.
# 2. Move the "_GoBack" bookmark from the very end of the document up to
#    the empty Heading-4 paragraph that sits right before the paragraph
#    holding the floating QR textbox.
# 3. Re-anchor/resize/lock the QR textbox ("Cuadro de texto 2") so it is
#    pinned to a fixed spot on the page instead of floating relative to
#    the paragraph.

$d = $word.ActiveDocument

# --- 1. Merge the placeholder runs -----------------------------------
$d.Content.Find.Execute("+++domicilio+++", $false, $false, $false, $false, `
    $false, $true, 1, $false, "+++domicilio+++", 2) | Out-Null

$d.Content.Find.Execute("+++localidad+++", $false, $false, $false, $false, `
    $false, $true, 1, $false, "+++localidad+++", 2) | Out-Null

# --- 2. Relocate the _GoBack bookmark ---------------------------------
# Paragraph 18 is the empty "Ttulo4" paragraph immediately before the
# paragraph that contains the anchored drawing (paragraph 19).
$anchorPara = $d.Paragraphs.Item(18)
$bmRange = $anchorPara.Range.Duplicate
$bmRange.Collapse(1)
$d.Bookmarks.Add("_GoBack", $bmRange) | Out-Null

# --- 3. Reposition / lock the QR textbox shape ------------------------
$qr = $d.Shapes.Item(3)
$qr.LockAnchor = $true
$qr.RelativeVerticalPosition = 1
$qr.Left = 346.9
$qr.Top = 554.8
$qr.Width = 103.18110236220473
$qr.Height = 94.67716535433071
